$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Neutrophils"-sender rows (old rows 6-9): the new
# TPM run no longer reports Neutrophils as a sending cluster, so those
# four rows disappear entirely and the sheet shrinks to A1:T5.
$ws.Rows("6:9").Delete()

# Row 2 (ECs -> ECs) keeps its cluster/gene labels; only the computed
# NATMI metrics change with the refreshed TPM values.
$ws.Range("G2").Value = 0.1806205
$ws.Range("H2").Value = 0.361241
$ws.Range("I2").Value = 0.7284509268949775
$ws.Range("J2").Value = 0.7284509268949775
$ws.Range("M2").Value = 1.058059
$ws.Range("N2").Value = 2.116118
$ws.Range("O2").Value = 0.806213901461467
$ws.Range("P2").Value = 0.806213901461467
$ws.Range("Q2").Value = 0.1911071456095
$ws.Range("R2").Value = 0.764428582438
$ws.Range("S2").Value = 0.5872872637952217
$ws.Range("T2").Value = 0.5872872637952217

# Row 3 (ECs -> MuSCs): target cluster relabelled from MuSCs (it stays
# MuSCs) but now sourced from the former "Neutrophils" slot in the
# refreshed data set, plus refreshed metric values.
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.1806205
$ws.Range("H3").Value = 0.361241
$ws.Range("I3").Value = 0.7284509268949775
$ws.Range("J3").Value = 0.7284509268949775
$ws.Range("M3").Value = 0.254321
$ws.Range("N3").Value = 0.508642
$ws.Range("O3").Value = 0.193786098538533
$ws.Range("P3").Value = 0.193786098538533
$ws.Range("Q3").Value = 0.0459355861805
$ws.Range("R3").Value = 0.183742344722
$ws.Range("S3").Value = 0.1411636630997558
$ws.Range("T3").Value = 0.1411636630997558

# Row 4 (was ECs -> Neutrophils, now MuSCs -> ECs): sending cluster
# switches to MuSCs and target cluster becomes ECs, with refreshed
# metric values.
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 0.067331
$ws.Range("H4").Value = 0.134662
$ws.Range("I4").Value = 0.2715490731050226
$ws.Range("J4").Value = 0.2715490731050226
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 1.058059
$ws.Range("N4").Value = 2.116118
$ws.Range("O4").Value = 0.806213901461467
$ws.Range("P4").Value = 0.806213901461467
$ws.Range("Q4").Value = 0.071240170529
$ws.Range("R4").Value = 0.284960682116
$ws.Range("S4").Value = 0.2189266376662454
$ws.Range("T4").Value = 0.2189266376662454

# Row 5 (was ECs -> Resolving-Mac, now MuSCs -> MuSCs): sending cluster
# switches to MuSCs and the obsolete "Resolving-Mac" target is replaced
# by MuSCs, with refreshed metric values.
$ws.Range("A5").Value = "MuSCs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.067331
$ws.Range("H5").Value = 0.134662
$ws.Range("I5").Value = 0.2715490731050226
$ws.Range("J5").Value = 0.2715490731050226
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.254321
$ws.Range("N5").Value = 0.508642
$ws.Range("O5").Value = 0.193786098538533
$ws.Range("P5").Value = 0.193786098538533
$ws.Range("Q5").Value = 0.017123687251
$ws.Range("R5").Value = 0.068494749004
$ws.Range("S5").Value = 0.05262243543877722
$ws.Range("T5").Value = 0.05262243543877722
